# Update cointegration results with refined numerical precision
# (values recomputed with higher-precision solver; comments added for clarity)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 0.6183949840649011
$ws.Range("E2").Value = 0.0007819938653903547
$ws.Range("F2").Value = 24.96723852706736
$ws.Range("I2").Value = 1.611382155183994
$ws.Range("J2").Value = '[ 0.31778081 -0.06268569]'
$ws.Range("D3").Value = 0.6325630073287644
$ws.Range("E3").Value = 0.01836905748085003
$ws.Range("F3").Value = 16.73684346747135
$ws.Range("I3").Value = 1.080193585219813
$ws.Range("J3").Value = '[0.17631894 0.04110196]'
$ws.Range("D4").Value = 0.7024686787693805
$ws.Range("E4").Value = 0.381460095988693
$ws.Range("F4").Value = 9.977877839429048
$ws.Range("I4").Value = 0.6439708692505662
$ws.Range("J4").Value = '[0.02421957 0.12041545]'
$ws.Range("D5").Value = 0.6035278543255299
$ws.Range("E5").Value = 0.0006166652298220817
$ws.Range("F5").Value = 20.22051492922491
$ws.Range("I5").Value = 1.305029264260077
$ws.Range("J5").Value = '[ 0.19538445 -0.25052375]'
$ws.Range("D6").Value = 0.6754184082610786
$ws.Range("E6").Value = 0.002798446101486107
$ws.Range("F6").Value = 18.89073830814724
$ws.Range("I6").Value = 1.219205663253405
$ws.Range("J6").Value = '[ 0.13661119 -0.35028069]'
$ws.Range("D7").Value = 0.8712663639086411
$ws.Range("E7").Value = 0.001267735632991838
$ws.Range("F7").Value = 18.6060452984982
$ws.Range("I7").Value = 1.200831615400386
$ws.Range("J7").Value = '[ 0.43588107 -0.54375418]'
$ws.Range("D8").Value = 0.8480505866909828
$ws.Range("E8").Value = 0.01197726617947072
$ws.Range("F8").Value = 16.00553066823936
$ws.Range("I8").Value = 1.032994757313293
$ws.Range("J8").Value = '[ 0.37952501 -0.46381186]'
$ws.Range("D9").Value = 0.7302709556458745
$ws.Range("E9").Value = 0.009545117298338521
$ws.Range("F9").Value = 15.67207148204413
$ws.Range("I9").Value = 1.011473347104686
$ws.Range("J9").Value = '[ 0.37999037 -0.04464737]'
$ws.Range("D10").Value = 0.7964449650503628
$ws.Range("E10").Value = 0.0234106584591638
$ws.Range("F10").Value = 10.45987912316099
$ws.Range("I10").Value = 0.6750791660908199
$ws.Range("J10").Value = '[ 0.40515237 -0.38743099]'
$ws.Range("D11").Value = 0.7386033235556743
$ws.Range("E11").Value = 0.01826057667263337
$ws.Range("F11").Value = 12.30911419801227
$ws.Range("I11").Value = 0.7944285445623402
$ws.Range("J11").Value = '[ 0.4154796  -0.03906522]'
$ws.Range("D12").Value = 0.7782429767873497
$ws.Range("E12").Value = 0.02953494372728297
$ws.Range("F12").Value = 10.78442102754372
$ws.Range("I12").Value = 0.6960250561525023
$ws.Range("J12").Value = '[ 0.33661067 -0.03340498]'
$ws.Range("D13").Value = 0.751536855936819
$ws.Range("E13").Value = 0.01912854279031435
$ws.Range("F13").Value = 18.3778024683985
$ws.Range("I13").Value = 1.186100854404426
$ws.Range("J13").Value = '[ 0.2126317  -0.31504718]'
$ws.Range("D14").Value = 0.7721853018678833
$ws.Range("E14").Value = 0.03016477551485748
$ws.Range("F14").Value = 9.293356106216711
$ws.Range("I14").Value = 0.5997919303367504
$ws.Range("J14").Value = '[ 0.26785419 -0.09700595]'
$ws.Range("D15").Value = 0.5505465875505141
$ws.Range("E15").Value = 0.0686855988563593
$ws.Range("F15").Value = 7.381506096692615
$ws.Range("I15").Value = 0.4764013925567863
$ws.Range("J15").Value = '[ 0.21649302 -0.02877683]'
$ws.Range("D16").Value = 0.6020398630407148
$ws.Range("E16").Value = 0.0631437496157131
$ws.Range("F16").Value = 8.524155079169752
$ws.Range("I16").Value = 0.5501478013959812
$ws.Range("J16").Value = '[ 0.22800415 -0.05137533]'
$ws.Range("D17").Value = 0.3725246422785028
$ws.Range("E17").Value = 0.08313837459882889
$ws.Range("F17").Value = 7.319725762113397
$ws.Range("I17").Value = 0.4724140982240822
$ws.Range("D18").Value = 0.7658830705050073
$ws.Range("E18").Value = 0.01911828951625773
$ws.Range("F18").Value = 15.74095588352838
$ws.Range("I18").Value = 1.015919136942513
$ws.Range("D19").Value = 0.6021386179596223
$ws.Range("E19").Value = 0.0008485689004768245
$ws.Range("F19").Value = 18.15947727430979
$ws.Range("I19").Value = 1.17201017627836
$ws.Range("J19").Value = '[ 0.55961217 -0.06370079]'
$ws.Range("D20").Value = 0.6099211632949364
$ws.Range("E20").Value = 0.002027452177843329
$ws.Range("F20").Value = 17.59091876843787
$ws.Range("I20").Value = 1.135315488175514
$ws.Range("J20").Value = '[ 0.47875303 -0.20064189]'
$ws.Range("D21").Value = 0.4275101408568411
$ws.Range("E21").Value = 0.03383156531121734
$ws.Range("F21").Value = 10.48519306123329
$ws.Range("I21").Value = 0.6767129241871717
$ws.Range("J21").Value = '[ 0.36061123 -0.50433333]'
$ws.Range("D22").Value = 0.5978225160495485
$ws.Range("E22").Value = 0.01126702359254859
$ws.Range("F22").Value = 12.38660744062619
$ws.Range("I22").Value = 0.799429947827665
$ws.Range("J22").Value = '[ 0.4710037  -0.18720619]'
$ws.Range("D23").Value = 0.6269288243491759
$ws.Range("E23").Value = 0.2750023459548397
$ws.Range("F23").Value = 7.097655564617231
$ws.Range("I23").Value = 0.4580817180909902
$ws.Range("J23").Value = '[ 0.12747955 -0.40589609]'
$ws.Range("D24").Value = 0.6420917582793985
$ws.Range("E24").Value = 0.1013521746878666
$ws.Range("F24").Value = 7.81891255800633
$ws.Range("I24").Value = 0.5046315456655887
$ws.Range("J24").Value = '[ 0.1408865  -0.12248264]'
$ws.Range("D25").Value = 0.6850898890537228
$ws.Range("E25").Value = 0.001150513855321031
$ws.Range("F25").Value = 16.78220018378531
$ws.Range("I25").Value = 1.083120901478951
$ws.Range("J25").Value = '[ 0.7632276  -0.22079367]'
$ws.Range("D26").Value = 0.703761072875736
$ws.Range("E26").Value = 0.002482864907776305
$ws.Range("F26").Value = 17.44921113920149
$ws.Range("I26").Value = 1.126169697191967
$ws.Range("J26").Value = '[ 0.60993769 -0.71030695]'
$ws.Range("D27").Value = 0.6463059403572415
$ws.Range("E27").Value = 0.003636998589662188
$ws.Range("F27").Value = 15.61973390592196
$ws.Range("I27").Value = 1.00809548710958
$ws.Range("J27").Value = '[ 0.26978576 -0.07183855]'
